$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(10013.64, 9928.26, 305.24, 307.86, $false, 0.86, 42613.765462962961, $true),
    @(9946.5499999999993, 10013.64, 307.68, 305.63, $false, -0.67, 42614.672662037039, $false),
    @(9883.89, 9946.5499999999993, 307.95999999999998, 306.02, $false, -0.63, 42615.7500462963, $false)
)

$startRow = 11
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
    $ws.Cells.Item($row, 8).Value = $rowData[7]

    # Copy the date format (style) from the row above so the new G cell
    # reuses the existing date-formatted style instead of creating a new one.
    $ws.Cells.Item($row - 1, 7).Copy() | Out-Null
    $ws.Cells.Item($row, 7).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# Column A width changed from 8.85546875 to 9 characters.
$ws.Columns.Item(1).ColumnWidth = 8.14
